$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates -------------------------------------------------
# Row 1 (A1:AN1): all cells go from 1 -> 21
$ws.Range("A1:AN1").Value = 21

# Rows 2-289: only column A (1 -> 6) and column AN (1 -> 10) change;
# the rest of the row (B:AM) stays as-is.
$ws.Range("A2:A289").Value = 6
$ws.Range("AN2:AN289").Value = 10

# Row 290 (A290:AN290): all cells go from 1 -> 2
$ws.Range("A290:AN290").Value = 2

# --- Conditional formatting ---------------------------------------------
# Swap the priority of the two existing "cellIs equal" rules so the
# "equals 2" (green) rule now ranks ahead of the "equals 1" (red) rule,
# and drop the duplicate red dxf that the original file carried.
$cfRange = $ws.Range("A1:AN300")
$fcOne = $cfRange.FormatConditions.Item(1)
$fcTwo = $cfRange.FormatConditions.Item(2)
$fcTwo.SetFirstPriority()

# --- Sheet view -----------------------------------------------------------
# Scroll the visible top-left cell and move the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 274
$win.ScrollColumn = 1
$ws.Range("AR280").Select()
